# Updated symbol list on Thu Jan 19 16:18:51 UTC 2023 with GitHub Actions
# Refreshes the Price (D), Volume(1h) (E) and Hora (G) columns for each
# coin row on Sheet1. Every one of these columns is stored as literal
# text in the workbook (e.g. "290.35", "-1.39%", "16"), so a leading
# apostrophe is used when writing through the Excel object model --
# otherwise Excel would silently reinterpret the digits/percentages as
# numbers instead of preserving the original text cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "'290.35"
    "E2" = "'-1.39%"
    "G2" = "'16"
    "D3" = "'30.69"
    "E3" = "'-4.15%"
    "G3" = "'16"
    "D4" = "'4.909"
    "E4" = "'-2.04%"
    "G4" = "'16"
    "D5" = "'0.07250"
    "E5" = "'-2.78%"
    "G5" = "'16"
    "D6" = "'2.370"
    "E6" = "'24.75%"
    "G6" = "'16"
    "D7" = "'7.654"
    "E7" = "'-1.29%"
    "G7" = "'16"
    "D8" = "'3.699"
    "E8" = "'-2.20%"
    "G8" = "'16"
    "D9" = "'0.8979"
    "E9" = "'-2.23%"
    "G9" = "'16"
    "D10" = "'0.1677"
    "E10" = "'-2.07%"
    "G10" = "'16"
    "D11" = "'0.08026"
    "E11" = "'3.31%"
    "G11" = "'16"
    "D12" = "'0.08137"
    "E12" = "'-2.02%"
    "G12" = "'16"
    "D13" = "'0.03065"
    "E13" = "'0.52%"
    "G13" = "'16"
    "D14" = "'0.1000"
    "E14" = "'0.05%"
    "G14" = "'16"
    "D15" = "'0.001495"
    "E15" = "'-0.97%"
    "G15" = "'16"
    "D16" = "'0.005727"
    "E16" = "'-2.13%"
    "G16" = "'16"
    "D17" = "'3.483"
    "E17" = "'0.56%"
    "G17" = "'16"
    "D18" = "'2.079"
    "E18" = "'-2.96%"
    "G18" = "'16"
    "D19" = "'0.3314"
    "E19" = "'-0.69%"
    "G19" = "'16"
    "D20" = "'0.1287"
    "E20" = "'-2.93%"
    "G20" = "'16"
    "D21" = "'3.961"
    "E21" = "'-9.91%"
    "G21" = "'16"
    "E22" = "'8.57%"
    "G22" = "'16"
    "D23" = "'0.04510"
    "G23" = "'16"
    "G24" = "'16"
    "D25" = "'0.004434"
    "E25" = "'9.64%"
    "G25" = "'16"
    "E26" = "'2.69%"
    "G26" = "'16"
    "D27" = "'0.0003390"
    "G27" = "'16"
    "G28" = "'16"
    "G29" = "'16"
    "G30" = "'16"
    "G31" = "'16"
    "G32" = "'16"
    "G33" = "'16"
    "G34" = "'16"
    "G35" = "'16"
    "G36" = "'16"
    "G37" = "'16"
    "G38" = "'16"
    "E39" = "'-2.72%"
    "G39" = "'16"
    "D40" = "'0.04351"
    "E40" = "'-3.10%"
    "G40" = "'16"
    "D41" = "'0.007308"
    "E41" = "'-3.04%"
    "G41" = "'16"
    "D42" = "'0.01002"
    "G42" = "'16"
    "D43" = "'0.1314"
    "E43" = "'-1.02%"
    "G43" = "'16"
    "D44" = "'0.002032"
    "E44" = "'-10.09%"
    "G44" = "'16"
    "D45" = "'0.009184"
    "E45" = "'-9.26%"
    "G45" = "'16"
    "D46" = "'0.00005707"
    "E46" = "'-6.75%"
    "G46" = "'16"
    "E47" = "'-1.22%"
    "G47" = "'16"
    "E48" = "'25.22%"
    "G48" = "'16"
    "D49" = "'0.002895"
    "E49" = "'-4.59%"
    "G49" = "'16"
    "E50" = "'-1.22%"
    "G50" = "'16"
    "E51" = "'-1.22%"
    "G51" = "'16"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
